$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28:97 down to 29:98
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly data point
$ws.Cells.Item(28, 1).Value = 9
$ws.Cells.Item(28, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 45014
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = 100112029
$ws.Cells.Item(28, 7).Value = "Orégano"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 16
$ws.Cells.Item(28, 11).Value = 18000
$ws.Cells.Item(28, 12).Value = 19000
$ws.Cells.Item(28, 13).Value = 18500
$ws.Cells.Item(28, 14).Value = "$/docena de atados"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 6167
$ws.Cells.Item(28, 17).Value = 3
$ws.Cells.Item(28, 18).Value = "Hortaliza"
